$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G5 / G18: store as literal text that looks like scientific notation
# (force text type via NumberFormat "@" so Excel doesn't re-interpret the
# large numeric-looking string as a real number, then drop the temporary
# number-format override so no stray cell style is left behind).
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5.22722e+16"
$ws.Range("G5").ClearFormats()

$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "2.72721e+16"
$ws.Range("G18").ClearFormats()

# I3: new tutor e-mail for JUAN MANUEL BENÍTEZ HERNÁNDEZ
$ws.Range("I3").Value = "Marchate1986@gmail.com"

# I4: updated tutor e-mail for JOSÉ DOMINGO CASTRO JUÁREZ
$ws.Range("I4").Value = "jdcastro@gmail.com"

# I20: new tutor e-mail for MARIA ISABEL OFICIAL SÁNCHEZ
$ws.Range("I20").Value = "oficialsanchezmz@gmail.com"
